# Add columns I (I0) and J (IF) with header labels and per-row values,
# matching style of the existing header row (bold, bordered, centered).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same style (s="1": bold, border, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the new header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data for columns I and J (rows 2-74).
$data = @(
    @(6,7),
    @(6,7),
    @(6,6),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(9,9),
    @(8,8),
    @(6,7),
    @(6,6),
    @(8,8),
    @(9,9),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(11,11),
    @(8,8),
    @(7,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,9),
    @(9,9),
    @(7,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,8),
    @(9,9),
    @(9,9),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,8),
    @(7,8),
    @(9,9),
    @(8,9),
    @(9,9),
    @(7,7),
    @(5,5),
    @(4,4),
    @(5,5),
    @(4,4),
    @(3,3)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
